$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save original row 2 values (A, D, G) before overwriting, since the
# rotation moves row2 -> row6, row3 -> row2, row6 -> row3.
$origA2 = $ws.Range("A2").Value()
$origD2 = $ws.Range("D2").Value()
$origG2 = $ws.Range("G2").Value()

$origA3 = $ws.Range("A3").Value()
$origD3 = $ws.Range("D3").Value()
$origG3 = $ws.Range("G3").Value()

$origA6 = $ws.Range("A6").Value()
$origD6 = $ws.Range("D6").Value()
$origG6 = $ws.Range("G6").Value()

# Row 2 <- old Row 3 (Instrument Data)
$ws.Range("A2").Value = $origA3
$ws.Range("D2").Value = $origD3
$ws.Range("G2").Value = $origG3

# Row 3 <- old Row 6 (Verification)
$ws.Range("A3").Value = $origA6
$ws.Range("D3").Value = $origD6
$ws.Range("G3").Value = $origG6

# Row 6 <- old Row 2 (Uncertainty)
$ws.Range("A6").Value = $origA2
$ws.Range("D6").Value = $origD2
$ws.Range("G6").Value = $origG2
